# edit.ps1 - applies the tracked changes to draft-gandhi-mpls-ioam-sr-03.pptx
#
#   1. Handout Master date placeholder: 11/9/20 -> 11/11/20
#   2. Slide 16 ("IOAM and Flow Indicator Label"): tighten paragraph line
#      spacing on the body placeholder from 21.2pt to 20.2pt (spcPts 2120 -> 2020)
#   3. Slide 5 ("Updates Since IETF-108"): reword a bullet:
#        "Need for different HbH and E2E Indicator Labels"
#        -> "Why different HbH and E2E Indicator Labels – Optimize processing on transit nodes"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Handout master date field: 11/9/20 -> 11/11/20 (best effort - some
#    hosts do not allow editing Handout Master / Notes Master placeholder
#    text; guarded so the rest of the script still runs either way).
# ---------------------------------------------------------------------------
try {
    $hf = $p.HandoutMaster.HeadersFooters.DateAndTime
    $hf.Value = "11/11/20"
} catch {
    Write-Output "HandoutMaster date edit not supported: $_"
}

# ---------------------------------------------------------------------------
# 2) Slide 16: line spacing 21.2 -> 20.2 on every paragraph of the body text
# ---------------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
for ($i = 1; $i -le $s16.Shapes.Count; $i++) {
    $shape = $s16.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($j = 1; $j -le $paraCount; $j++) {
            $para = $tr.Paragraphs($j)
            if ([math]::Round($para.ParagraphFormat.SpaceWithin, 1) -eq 21.2) {
                $para.ParagraphFormat.SpaceWithin = 20.2
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 5: reword bullet, preserving the existing run split / formatting
#    (the middle run "HbH" keeps its err="1" flag untouched).
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$content = $s5.Shapes.Item(2)
$tr5 = $content.TextFrame.TextRange

$full = $tr5.Text
$old1 = "Need for different "
$idx1 = $full.IndexOf($old1)
if ($idx1 -ge 0) {
    $tr5.Characters($idx1 + 1, $old1.Length).Text = "Why different "
}

$full = $tr5.Text
$old2 = " and E2E Indicator Labels"
$idx2 = $full.IndexOf($old2)
if ($idx2 -ge 0) {
    $dash = [char]0x2013
    $tr5.Characters($idx2 + 1, $old2.Length).Text = " and E2E Indicator Labels " + $dash + " Optimize processing on transit nodes"
}
